$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 310, pushing the existing rows 310-329
# down to 312-331 (so the former rows 328/329 become the new rows 330/331).
$ws.Range("A310:A311").EntireRow.Insert()

# New weekly entries for the two inserted rows.
$ws.Cells.Item(310, 1).Value2 = 11
$ws.Cells.Item(310, 2).Value2 = 'Vega Monumental Concepción'
$ws.Cells.Item(310, 3).Value2 = 'Bíobío'
$ws.Cells.Item(310, 4).Value2 = 44826
$ws.Cells.Item(310, 5).Value2 = 8
$ws.Cells.Item(310, 6).Value2 = 100114001
$ws.Cells.Item(310, 7).Value2 = 'Papa'
$ws.Cells.Item(310, 8).Value2 = 'Asterix'
$ws.Cells.Item(310, 9).Value2 = '1a (guarda)'
$ws.Cells.Item(310, 10).Value2 = 4000
$ws.Cells.Item(310, 11).Value2 = 7500
$ws.Cells.Item(310, 12).Value2 = 8000
$ws.Cells.Item(310, 13).Value2 = 7750
$ws.Cells.Item(310, 14).Value2 = '$/saco 25 kilos'
$ws.Cells.Item(310, 15).Value2 = 'Región de La Araucanía'
$ws.Cells.Item(310, 16).Value2 = 310
$ws.Cells.Item(310, 17).Value2 = 25
$ws.Cells.Item(310, 18).Value2 = 'Hortaliza'

$ws.Cells.Item(311, 1).Value2 = 11
$ws.Cells.Item(311, 2).Value2 = 'Vega Monumental Concepción'
$ws.Cells.Item(311, 3).Value2 = 'Bíobío'
$ws.Cells.Item(311, 4).Value2 = 44826
$ws.Cells.Item(311, 5).Value2 = 8
$ws.Cells.Item(311, 6).Value2 = 100114001
$ws.Cells.Item(311, 7).Value2 = 'Papa'
$ws.Cells.Item(311, 8).Value2 = 'Patagonia'
$ws.Cells.Item(311, 9).Value2 = '1a (guarda)'
$ws.Cells.Item(311, 10).Value2 = 4000
$ws.Cells.Item(311, 11).Value2 = 7500
$ws.Cells.Item(311, 12).Value2 = 8000
$ws.Cells.Item(311, 13).Value2 = 7750
$ws.Cells.Item(311, 14).Value2 = '$/saco 25 kilos'
$ws.Cells.Item(311, 15).Value2 = 'Región de La Araucanía'
$ws.Cells.Item(311, 16).Value2 = 310
$ws.Cells.Item(311, 17).Value2 = 25
$ws.Cells.Item(311, 18).Value2 = 'Hortaliza'
